$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $savedStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $savedStyle
}

Set-TextValue 'D2' '28.012.32'
Set-TextValue 'E2' '  -0.30%  '
Set-TextValue 'D3' '1.764.25'
Set-TextValue 'E3' '  -2.87%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.26%  '
Set-TextValue 'D5' '339.19'
Set-TextValue 'E5' '  +0.32%  '
Set-TextValue 'D6' '0.9976'
Set-TextValue 'E7' '  -4.10%  '
Set-TextValue 'D8' '0.3362'
Set-TextValue 'E8' '  -3.64%  '
Set-TextValue 'D9' '45.83'
Set-TextValue 'E9' '  -5.43%  '
Set-TextValue 'D10' '1.133'
Set-TextValue 'E10' '  -5.66%  '
Set-TextValue 'D11' '0.07201'
Set-TextValue 'E11' '  -5.15%  '
Set-TextValue 'D12' '22.81'
Set-TextValue 'E12' '  +2.70%  '
Set-TextValue 'D13' '0.9982'
Set-TextValue 'E13' '  +0.20%  '
Set-TextValue 'D14' '6.215'
Set-TextValue 'E14' '  -4.97%  '
Set-TextValue 'D15' '7.207'
Set-TextValue 'E15' '  +0.15%  '
Set-TextValue 'D16' '1.761.70'
Set-TextValue 'E16' '  -2.82%  '
Set-TextValue 'D17' '0.00001055'
Set-TextValue 'E17' '  -4.74%  '
Set-TextValue 'D18' '0.06584'
Set-TextValue 'E18' '  -1.87%  '
Set-TextValue 'D19' '80.67'
Set-TextValue 'E19' '  -5.50%  '
Set-TextValue 'D20' '0.9983'
Set-TextValue 'E20' '  +0.12%  '
Set-TextValue 'D21' '17.01'
Set-TextValue 'E22' '  -4.40%  '
Set-TextValue 'D23' '27.989.47'
Set-TextValue 'E23' '  -0.37%  '
Set-TextValue 'D24' '11.77'
Set-TextValue 'E24' '  -8.32%  '
Set-TextValue 'D25' '2.378'
Set-TextValue 'E25' '  -1.20%  '
Set-TextValue 'D26' '153.31'
Set-TextValue 'E26' '  -0.93%  '
Set-TextValue 'D27' '2.357'
Set-TextValue 'E27' '  -8.31%  '
Set-TextValue 'D28' '19.77'
Set-TextValue 'E28' '  -7.60%  '
Set-TextValue 'B29' 'ImmutableX'
Set-TextValue 'C29' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D29' '1.295'
Set-TextValue 'E29' '  -14.55%  '
Set-TextValue 'B30' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C30' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D30' '1.961.47'
Set-TextValue 'E30' '  -3.11%  '
Set-TextValue 'D31' '130.89'
Set-TextValue 'E31' '  -3.48%  '
Set-TextValue 'D32' '4.017'
Set-TextValue 'E32' '  -0.52%  '
Set-TextValue 'D33' '5.853'
Set-TextValue 'E33' '  -4.74%  '
Set-TextValue 'D34' '0.08750'
Set-TextValue 'E34' '  -1.11%  '
Set-TextValue 'D35' '12.33'
Set-TextValue 'E35' '  -7.35%  '
Set-TextValue 'E36' '  -3.23%  '
Set-TextValue 'D37' '0.6612'
Set-TextValue 'E37' '  -5.04%  '
Set-TextValue 'D38' '0.06215'
Set-TextValue 'E38' '  -5.27%  '
Set-TextValue 'D39' '5.159'
Set-TextValue 'E39' '  -6.70%  '
Set-TextValue 'D40' '0.2111'
Set-TextValue 'E40' '  -5.36%  '
Set-TextValue 'D41' '1.214'
Set-TextValue 'E41' '  -4.22%  '
Set-TextValue 'D42' '1.448'
Set-TextValue 'E42' '  -10.44%  '
Set-TextValue 'D43' '8.048'
Set-TextValue 'E43' '  -5.63%  '
Set-TextValue 'D44' '0.9973'
Set-TextValue 'E44' '  +0.17%  '
Set-TextValue 'D45' '13.67'
Set-TextValue 'E45' '  -6.76%  '
Set-TextValue 'D46' '3.836'
Set-TextValue 'E46' '  -1.04%  '
Set-TextValue 'D47' '0.6042'
Set-TextValue 'E47' '  -7.38%  '
Set-TextValue 'D48' '130.21'
Set-TextValue 'E48' '  -1.80%  '
Set-TextValue 'D49' '2.017'
Set-TextValue 'E49' '  -7.08%  '
Set-TextValue 'D50' '0.07239'
Set-TextValue 'E50' '  +0.30%  '
Set-TextValue 'D51' '1.178'
Set-TextValue 'E51' '  +1.45%  '
